# Fruta / hortaliza, semanal
# Insert a new weekly record as row 3, shifting existing rows (old rows 3..19) down to rows 4..20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 3 (pushes current rows 3-19 down to 4-20,
# and copies formatting - incl. the date number format on column D - from row 2 above).
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the new weekly record.
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(3, 3).Value = "Ñuble"
$ws.Cells.Item(3, 4).Value = 45222
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = "Fruta"
$ws.Cells.Item(3, 7).Value = 100102
$ws.Cells.Item(3, 8).Value = "Cítricos"
$ws.Cells.Item(3, 9).Value = 100102006
$ws.Cells.Item(3, 10).Value = "Pomelo"
$ws.Cells.Item(3, 11).Value = "Start Ruby"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 80
$ws.Cells.Item(3, 14).Value = 15000
$ws.Cells.Item(3, 15).Value = 15000
$ws.Cells.Item(3, 16).Value = 15000
$ws.Cells.Item(3, 17).Value = "`$/caja 14 kilos granel"
$ws.Cells.Item(3, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(3, 19).Value = 1071
$ws.Cells.Item(3, 20).Value = 14
